$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.573811
$ws.Range("H2").Value = 13.721433
$ws.Range("I2").Value = 0.1659009079913533
$ws.Range("J2").Value = 0.1659009079913533
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 9.271642409868665
$ws.Range("R2").Value = 83.44478168881801
$ws.Range("S2").Value = 0.001094329598772298
$ws.Range("T2").Value = 0.001094329598772298
$ws.Range("G3").Value = 4.573811
$ws.Range("H3").Value = 13.721433
$ws.Range("I3").Value = 0.1659009079913533
$ws.Range("J3").Value = 0.1659009079913533
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 1172.927784620959
$ws.Range("R3").Value = 10556.35006158863
$ws.Range("S3").Value = 0.1384403685119383
$ws.Range("T3").Value = 0.1384403685119382
$ws.Range("G4").Value = 4.573811
$ws.Range("H4").Value = 13.721433
$ws.Range("I4").Value = 0.1659009079913533
$ws.Range("J4").Value = 0.1659009079913533
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 223.3861443491224
$ws.Range("R4").Value = 2010.475299142101
$ws.Range("S4").Value = 0.02636620988064276
$ws.Range("T4").Value = 0.02636620988064275
$ws.Range("I5").Value = 0.5322852674812913
$ws.Range("J5").Value = 0.5322852674812913
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 29.74762899058422
$ws.Range("R5").Value = 267.728660915258
$ws.Range("S5").Value = 0.003511105094286561
$ws.Range("T5").Value = 0.003511105094286561
$ws.Range("I6").Value = 0.5322852674812913
$ws.Range("J6").Value = 0.5322852674812913
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.4441794169530783
$ws.Range("T6").Value = 0.4441794169530783
$ws.Range("I7").Value = 0.5322852674812913
$ws.Range("J7").Value = 0.5322852674812913
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 716.7239470605202
$ws.Range("R7").Value = 6450.515523544682
$ws.Range("S7").Value = 0.08459474543392652
$ws.Range("T7").Value = 0.08459474543392648
$ws.Range("G8").Value = 8.320867
$ws.Range("H8").Value = 24.962601
$ws.Range("I8").Value = 0.3018138245273554
$ws.Range("J8").Value = 0.3018138245273554
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 16.86735708232733
$ws.Range("R8").Value = 151.806213740946
$ws.Range("S8").Value = 0.001990849872359757
$ws.Range("T8").Value = 0.001990849872359756
$ws.Range("G9").Value = 8.320867
$ws.Range("H9").Value = 24.962601
$ws.Range("I9").Value = 0.3018138245273554
$ws.Range("J9").Value = 0.3018138245273554
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 2133.838957586057
$ws.Range("R9").Value = 19204.55061827452
$ws.Range("S9").Value = 0.2518564701993209
$ws.Range("T9").Value = 0.2518564701993209
$ws.Range("G10").Value = 8.320867
$ws.Range("H10").Value = 24.962601
$ws.Range("I10").Value = 0.3018138245273554
$ws.Range("J10").Value = 0.3018138245273554
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 406.3933548569997
$ws.Range("R10").Value = 3657.540193712997
$ws.Range("S10").Value = 0.04796650445567476
$ws.Range("T10").Value = 0.04796650445567475
